$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 261 (shifts old rows 261-382 down to 263-384)
$ws.Rows("261:262").Insert()

# New row 261: same as old row 261 but with updated D/J/K/L/M/P values
$ws.Range("A261").Value = 6
$ws.Range("B261").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C261").Value = "Metropolitana"
$ws.Range("D261").Value = 44466
$ws.Range("E261").Value = 13
$ws.Range("F261").Value = 100112017
$ws.Range("G261").Value = "Apio"
$ws.Range("H261").Value = "Americana (o)"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 2000
$ws.Range("K261").Value = 6000
$ws.Range("L261").Value = 7000
$ws.Range("M261").Value = 6600
$ws.Range("N261").Value = "`$/docena de matas"
$ws.Range("O261").Value = "Región de Coquimbo"
$ws.Range("P261").Value = 1100
$ws.Range("Q261").Value = 6
$ws.Range("R261").Value = "Hortaliza"

# New row 262: same as old row 262 but with updated D/K/L/M/P values
$ws.Range("A262").Value = 6
$ws.Range("B262").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C262").Value = "Metropolitana"
$ws.Range("D262").Value = 44466
$ws.Range("E262").Value = 13
$ws.Range("F262").Value = 100112017
$ws.Range("G262").Value = "Apio"
$ws.Range("H262").Value = "Americana (o)"
$ws.Range("I262").Value = "Segunda"
$ws.Range("J262").Value = 600
$ws.Range("K262").Value = 5000
$ws.Range("L262").Value = 5000
$ws.Range("M262").Value = 5000
$ws.Range("N262").Value = "`$/docena de matas"
$ws.Range("O262").Value = "Región de Coquimbo"
$ws.Range("P262").Value = 833
$ws.Range("Q262").Value = 6
$ws.Range("R262").Value = "Hortaliza"
